$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Set column C (Public) to TRUE for rows 15 through 32
$ws.Range("C15:C32").Value = $true

# Update the selected cell/range as recorded in the workbook view
$ws.Range("C15:C38").Select()
